# Update "想去人数" (interested-people count) figures to the latest scrape,
# mirroring the same four events on both the "展览" sheet and the combined
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (rows 3-6 hold the four events)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F3").Value = 183
$wsExhibit.Range("F4").Value = 778
$wsExhibit.Range("F5").Value = 68
$wsExhibit.Range("F6").Value = 3

# Sheet 4: 全部类型 (same four events, shifted down one row)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F4").Value = 183
$wsAll.Range("F5").Value = 778
$wsAll.Range("F6").Value = 68
$wsAll.Range("F7").Value = 3
